$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 8 obsolete data rows (old A460..A467) so the remaining rows
# shift up, re-aligning tapeid A468.. with rows 44..62 as in the new data.
$ws.Range("A44:M51").EntireRow.Delete()
